$d = $word.ActiveDocument

# The document's headers/footers each contain an inline picture whose
# wp:docPr/@name and pic:cNvPr/@name were swapped:
#   - the two Pearson logo pictures (footer1.xml / footer2.xml) were
#     named "image2.png" and should become "image1.png"
#   - the BTEC logo picture (header1.xml) was named "image1.jpg" and
#     should become "image2.jpg"
# These are drawing-object display names (not the underlying media part
# file names, which are untouched), so they are edited by round-tripping
# the document's WordOpenXML and renaming the relevant name="..." values.

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml
